# Bangphancong.xlsx update — "update bang phan cong TrThinh"
#
# Adds the new assignee (Nguyễn Trường Thịnh, MSSV 15520844) as row 5 of the
# "Bảng phân công" sheet, selects the newly entered name cell (mirroring the
# author leaving the cursor on B5 after data entry), and updates the sheet's
# base/Normal font from Calibri to Arial.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New assignment row: MSSV in column A (numeric), full name in column B.
$ws.Range("A5").Value = 15520844
$ws.Range("B5").Value = "Nguyễn Trường Thịnh"

# Leave the selection on the cell that was just filled in, like the author did.
$ws.Range("B5").Select() | Out-Null

# Switch the workbook's base font from Calibri to Arial.
$wb.Styles.Item("Normal").Font.Name = "Arial"
